# Update "想去人数" (want-to-go count) values in column F
# on the "展览" and "全部类型" sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 93
$ws1.Range("F5").Value  = 48
$ws1.Range("F6").Value  = 29
$ws1.Range("F8").Value  = 54
$ws1.Range("F9").Value  = 8375
$ws1.Range("F10").Value = 779
$ws1.Range("F11").Value = 307
$ws1.Range("F12").Value = 1122
$ws1.Range("F13").Value = 887
$ws1.Range("F14").Value = 70
$ws1.Range("F15").Value = 43
$ws1.Range("F17").Value = 152
$ws1.Range("F19").Value = 220
$ws1.Range("F20").Value = 926

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 93
$ws4.Range("F6").Value  = 48
$ws4.Range("F7").Value  = 29
$ws4.Range("F10").Value = 54
$ws4.Range("F11").Value = 8375
$ws4.Range("F12").Value = 779
$ws4.Range("F13").Value = 307
$ws4.Range("F14").Value = 1122
$ws4.Range("F15").Value = 887
$ws4.Range("F16").Value = 70
$ws4.Range("F17").Value = 43
$ws4.Range("F19").Value = 152
$ws4.Range("F21").Value = 220
$ws4.Range("F22").Value = 926
